$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the sequential numbering values added in column A for rows 26, 28, 30, 31
$ws.Range("A26").Value = 1
$ws.Range("A28").Value = 2
$ws.Range("A30").Value = 3
$ws.Range("A31").Value = 4

# Update the selected cell/range shown in the sheet view
$ws.Range("O29").Select()
